$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "Amount" column (C) one column to the right (D),
# carrying its values, number formats and styling with it.
$ws.Range("C1:C3").Copy($ws.Range("D1:D3"))

# Clear the old column C so it can become the new "Budget Version" column
# with its own (default) formatting instead of the Amount column's style.
$ws.Range("C1:C3").ClearFormats()

# Populate the new "Budget Version" column.
$ws.Range("C1").Value = "Budget Version"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1

# Restore the active selection to C3, matching the saved view state.
$ws.Range("C3").Select()
